# Angular JS Definitions workbook - add a new "AngularJS: The Big Picture"
# section (header row + two sub-items) to the end of the Definitions sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Definitions")

# --- Row 54: new section header, merged across A:B like the "Angular JS"
#     header in row 2 (style carried over via PasteSpecial formats) ---
$ws.Range("A54").Value = "AngularJS: The Big Picture"
$ws.Range("A2:B2").Copy()
$ws.Range("A54:B54").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A54:B54").Merge()

# --- Row 55: sub-item, formatted like the existing "Controller capabilities"
#     label in A31 (plain single-column label, no merge) ---
$ws.Range("A55").Value = "Angular JS Benefits"
$ws.Range("A31").Copy()
$ws.Range("A55").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 56: another sub-item under "Angular JS Benefits" ---
$ws.Range("A56").Value = "Code Reduction"
$ws.Range("A31").Copy()
$ws.Range("A56").PasteSpecial(-4122)   # xlPasteFormats

# Clear the clipboard/marching-ants state and leave the selection where a
# user would end up after typing the last new cell.
$excel.CutCopyMode = $false
$ws.Range("B56").Select()
